$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.571.19"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "2.077.20"
$ws.Range("E3").Value = "  +4.03%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").Value = "'236.31"
$ws.Range("E5").Value = "  -3.05%  "
$ws.Range("E6").Value = "  +1.90%  "
$ws.Range("D7").Value = "'58.46"
$ws.Range("E7").Value = "  +5.88%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.384"
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("D10").Value = "'58.25"
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("D11").Value = "'0.0764"
$ws.Range("E11").Value = "  +0.86%  "
$ws.Range("D12").Value = "'0.102"
$ws.Range("E12").Value = "  +3.63%  "
$ws.Range("D13").Value = "2.382.08"
$ws.Range("E13").Value = "  +3.88%  "
$ws.Range("D14").Value = "'14.57"
$ws.Range("E14").Value = "  +3.02%  "
$ws.Range("D15").Value = "'21.15"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").Value = "'0.779"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").Value = "'5.26"
$ws.Range("E17").Value = "  +4.04%  "
$ws.Range("D18").Value = "2.153.29"
$ws.Range("E18").Value = "  +7.72%  "
$ws.Range("D19").Value = "37.672.89"
$ws.Range("E19").Value = "  +1.82%  "
$ws.Range("D20").Value = "'6.03"
$ws.Range("E20").Value = "  +19.71%  "
$ws.Range("D21").Value = "'68.50"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "0.0₃0815"
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("D23").Value = "'224.26"
$ws.Range("E23").Value = "  -2.03%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  +3.08%  "
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("D27").Value = "'162.94"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").Value = "'8.90"
$ws.Range("E28").Value = "  +2.45%  "
$ws.Range("E29").Value = "  +5.65%  "
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").Value = "'1.39"
$ws.Range("E30").Value = "  +5.89%  "
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'19.33"
$ws.Range("E31").Value = "  +0.71%  "
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").Value = "'4.50"
$ws.Range("E33").Value = "  +1.24%  "
$ws.Range("D34").Value = "'0.0628"
$ws.Range("E34").Value = "  +2.50%  "
$ws.Range("D35").Value = "'2.58"
$ws.Range("E35").Value = "  +10.50%  "
$ws.Range("D36").Value = "'4.41"
$ws.Range("E36").Value = "  +3.74%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").Value = "'3.38"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("B38").Value = "BinanceUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").Value = "'5.96"
$ws.Range("E39").Value = "  +13.75%  "
$ws.Range("E40").Value = "  -0.68%  "
$ws.Range("E41").Value = "  -4.80%  "
$ws.Range("D42").Value = "'0.0967"
$ws.Range("E42").Value = "  +9.48%  "
$ws.Range("D43").Value = "1.478.18"
$ws.Range("E43").Value = "  +2.59%  "
$ws.Range("D44").Value = "'4.34"
$ws.Range("E44").Value = "  +23.35%  "
$ws.Range("D45").Value = "'95.40"
$ws.Range("E45").Value = "  +7.60%  "
$ws.Range("D46").Value = "'16.54"
$ws.Range("E46").Value = "  +7.88%  "
$ws.Range("D47").Value = "'0.0211"
$ws.Range("E47").Value = "  +3.18%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("E49").Value = "  +2.42%  "
$ws.Range("D50").Value = "'7.34"
$ws.Range("E50").Value = "  +9.39%  "
$ws.Range("E51").Value = "  +1.53%  "
